# Update (Analyze PO & Forecast)
$wb = $excel.ActiveWorkbook

# --- Sheet "Forecast Comparison" ---
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsForecast.Range("D12").Value = 3
$wsForecast.Range("D17").Value = 3

# --- Sheet "Summary" ---
# Force text format so these remain stored as strings (matching the
# original inlineStr cell type) instead of Excel auto-coercing the
# digit-only / date-like text into a number / date serial.
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B9").NumberFormat = "@"
$wsSummary.Range("B9").Value = "55"

$wsSummary.Range("B13").NumberFormat = "@"
$wsSummary.Range("B13").Value = "2025-02-16"
